$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.962.08"
$ws.Range("E2").Value = "  -1.60%  "

$ws.Range("D3").Value = "1.905.71"
$ws.Range("E3").Value = "  -3.12%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4592"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3825"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07710"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9795"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.03"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.01%  "

$ws.Range("D12").Value = "1.922.15"
$ws.Range("E12").Value = "  -3.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.935"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.660"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07036"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.73"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009488"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.63"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").Value = "28.957.12"
$ws.Range("E21").Value = "  -1.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.302"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.18%  "

$ws.Range("E23").Value = "  -2.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.093"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.08"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.644"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.35"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.852"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09263"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8663"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.077"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.247"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.952"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.72%  "

$ws.Range("E35").Value = "  -1.66%  "

$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02037"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5499"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.402"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.22%  "

$ws.Range("E41").Value = "  -2.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.312"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.778"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5175"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.31"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06833"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.061"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000002585"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -17.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.61"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.774"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.09%  "
